$wb = $excel.ActiveWorkbook

# ---------- Typography sheet: new row 8 (ScrollWhellTxt / ariblk.ttf) ----------
$wsTypo = $wb.Worksheets.Item("Typography")
$wsTypo.Range("B7:J7").Copy($wsTypo.Range("B8:J8"))
$wsTypo.Range("B8").Value = "ScrollWhellTxt"
$wsTypo.Range("C8").Value = "ariblk.ttf"
$wsTypo.Range("D8").Value = 50
$wsTypo.Range("E8").Value = 4
$wsTypo.Range("H8").Value = ""
$wsTypo.Range("I8").Value = ""

# ---------- Translation sheet: BME680 sensor texts ----------
$wsTr = $wb.Worksheets.Item("Translation")

# Existing texts whose wording changed (add units/spacing)
$wsTr.Range("F6").Value = "Temperature: <value>°C"
$wsTr.Range("F21").Value = "Humidity: <value>%"

# Rows with a unique "Default / Left / LTR" template + custom text (based on row 7)
$wsTr.Range("C7:E7").Copy($wsTr.Range("C52:E52"))
$wsTr.Range("B52").Value = "SingleUseId66"
$wsTr.Range("F52").Value = "Menu Element: <value>"
$wsTr.Range("C7:E7").Copy($wsTr.Range("C54:E54"))
$wsTr.Range("B54").Value = "SingleUseId68"
$wsTr.Range("F54").Value = "Pressure: <value>mmHg"
$wsTr.Range("C7:E7").Copy($wsTr.Range("C67:E67"))
$wsTr.Range("B67").Value = "SingleUseId81"
$wsTr.Range("F67").Value = "CO2 Concentration: <value>ppm"

# Rows with "Default / Left / LTR / 0" template (based on row 7)
$wsTr.Range("B7:F7").Copy($wsTr.Range("B53:F53"))
$wsTr.Range("B53").Value = "SingleUseId67"
$wsTr.Range("B7:F7").Copy($wsTr.Range("B55:F55"))
$wsTr.Range("B55").Value = "SingleUseId69"
$wsTr.Range("B7:F7").Copy($wsTr.Range("B68:F68"))
$wsTr.Range("B68").Value = "SingleUseId82"

# Rows with "Default / Right / LTR / <value>" template (based on row 9)
$wsTr.Range("B9:F9").Copy($wsTr.Range("B56:F56"))
$wsTr.Range("B56").Value = "SingleUseId70"
$wsTr.Range("B9:F9").Copy($wsTr.Range("B57:F57"))
$wsTr.Range("B57").Value = "SingleUseId71"
$wsTr.Range("B9:F9").Copy($wsTr.Range("B58:F58"))
$wsTr.Range("B58").Value = "SingleUseId72"
$wsTr.Range("B9:F9").Copy($wsTr.Range("B59:F59"))
$wsTr.Range("B59").Value = "SingleUseId73"
$wsTr.Range("B9:F9").Copy($wsTr.Range("B60:F60"))
$wsTr.Range("B60").Value = "SingleUseId74"
$wsTr.Range("B9:F9").Copy($wsTr.Range("B61:F61"))
$wsTr.Range("B61").Value = "SingleUseId75"
$wsTr.Range("B9:F9").Copy($wsTr.Range("B62:F62"))
$wsTr.Range("B62").Value = "SingleUseId76"
$wsTr.Range("B9:F9").Copy($wsTr.Range("B63:F63"))
$wsTr.Range("B63").Value = "SingleUseId77"
$wsTr.Range("B9:F9").Copy($wsTr.Range("B64:F64"))
$wsTr.Range("B64").Value = "SingleUseId78"
$wsTr.Range("B9:F9").Copy($wsTr.Range("B65:F65"))
$wsTr.Range("B65").Value = "SingleUseId79"
$wsTr.Range("B9:F9").Copy($wsTr.Range("B66:F66"))
$wsTr.Range("B66").Value = "SingleUseId80"
$wsTr.Range("B9:F9").Copy($wsTr.Range("B69:F69"))
$wsTr.Range("B69").Value = "SingleUseId83"
$wsTr.Range("B9:F9").Copy($wsTr.Range("B70:F70"))
$wsTr.Range("B70").Value = "SingleUseId84"
$wsTr.Range("B9:F9").Copy($wsTr.Range("B71:F71"))
$wsTr.Range("B71").Value = "SingleUseId85"
$wsTr.Range("B9:F9").Copy($wsTr.Range("B72:F72"))
$wsTr.Range("B72").Value = "SingleUseId86"
$wsTr.Range("B9:F9").Copy($wsTr.Range("B73:F73"))
$wsTr.Range("B73").Value = "SingleUseId87"
$wsTr.Range("B9:F9").Copy($wsTr.Range("B74:F74"))
$wsTr.Range("B74").Value = "SingleUseId88"
$wsTr.Range("B9:F9").Copy($wsTr.Range("B75:F75"))
$wsTr.Range("B75").Value = "SingleUseId89"
$wsTr.Range("B9:F9").Copy($wsTr.Range("B76:F76"))
$wsTr.Range("B76").Value = "SingleUseId90"
$wsTr.Range("B9:F9").Copy($wsTr.Range("B77:F77"))
$wsTr.Range("B77").Value = "SingleUseId91"
$wsTr.Range("B9:F9").Copy($wsTr.Range("B78:F78"))
$wsTr.Range("B78").Value = "SingleUseId92"
$wsTr.Range("B9:F9").Copy($wsTr.Range("B79:F79"))
$wsTr.Range("B79").Value = "SingleUseId93"

Write-Host "BME680 sensor texts added"
